$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (62) to the gebiedsdefinities table for the "Niet te
# lokaliseren" area (elz93), continuing the volgnr sequence (volgnr 61).
$ws.Range("A62").Value = 61
$ws.Range("B62").Value = "elz93"
$ws.Range("C62").Value = "Niet te lokaliseren"
$ws.Range("D62").Value = "Niet te lokaliseren"

# Match the numeric formatting already used by the rest of column A
# (integer number format) so the new cell keeps the same style as A61.
$ws.Range("A62").NumberFormat = $ws.Range("A61").NumberFormat
